# Reformat the "Stream apportionment" template:
#  - split the old "label: {value}" single-cell rows into a bold label
#    cell plus a separate value cell
#  - rename a couple of labels ("Generated:" -> "Date generated:",
#    "Streams" header -> "GNIS Name")
#  - turn the stream rows into a real Excel Table (ListObject)
#  - widen / re-size the data columns to fit the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: "Date generated:" label + value -------------------------------
$ws.Range("A2").Value = "Date generated: "
$ws.Range("A2").Font.Bold = $true
$ws.Range("B2").Value = "{d.generated}"

# --- Row 3: "Point of interest:" label + value -----------------------------
$ws.Range("A3").Value = "Point of interest:"
$ws.Range("A3").Font.Bold = $true
$ws.Range("B3").Value = "{d.point}"

# --- Row 4: "Weighting factor:" label + value ------------------------------
$ws.Range("A4").Value = "Weighting factor:"
$ws.Range("A4").Font.Bold = $true
$ws.Range("B4").Value = "{d.weighting_factor}"

# --- Row 5: blank bold spacer row ------------------------------------------
$ws.Range("A5").Value = " "
$ws.Range("A5").Font.Bold = $true

# clear anything left over from the old B5/C5 and old row 9 content
$ws.Range("B5:D5").ClearContents()
$ws.Range("A9:D9").ClearContents()

# --- Row 6: table header row ------------------------------------------------
$ws.Range("A6").Value = "GNIS Name"
$ws.Range("B6").Value = "Distance (m)"
$ws.Range("C6").Value = "Apportionment (%)"

# --- Row 7: first (templated) data row --------------------------------------
$ws.Range("A7").Value = "{d.streams[i].gnis_name}"
$ws.Range("B7").Value = "{d.streams[i].distance}"
$ws.Range("C7").Value = "{d.streams[i].apportionment}"

# --- Row 8: repeat marker row ------------------------------------------------
$ws.Range("A8").Value = "{d.streams[i+1].gnis_name}"
$ws.Range("B8").ClearContents()
$ws.Range("C8").ClearContents()

# Turn A6:C8 into a proper Excel table (ListObject) with banded rows.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A6:C8"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleMedium2"

# --- Column widths for the new layout ---------------------------------------
$ws.Columns("A").ColumnWidth = 18.416666666666668
$ws.Columns("B").ColumnWidth = 11.916666666666666
$ws.Columns("C").ColumnWidth = 17.583333333333332

# --- Selection / view tidy-up ------------------------------------------------
[void]$ws.Range("B15").Select()
